$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J15").Value = 116
$ws.Range("H16").Value = 140
$ws.Range("I16").Value = 133
$ws.Range("I17").Value = 140
$ws.Range("J17").Value = 110
$ws.Range("J18").Value = 112

$ws.Range("J16").Select()
